$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(" Dubai (DSC)", " October 13 2020", "Super Kings won by 20 runs", "Sunrisers Hyderabad", "Chennai Super Kings", "Rashid Khan ", "14", "8", "1", "1", "175.00"),
    @(" Abu Dhabi", " October 18 2020", "Match tied (KKR won the one-over eliminator)", "Sunrisers Hyderabad", "Kolkata Knight Riders", "Rashid Khan ", "1", "2", "0", "0", "50.00"),
    @(" Sharjah", " October 04 2020", "Mumbai won by 34 runs", "Sunrisers Hyderabad", "Mumbai Indians", "Rashid Khan ", "3", "7", "0", "0", "42.85"),
    @(" Dubai (DSC)", " September 21 2020", "RCB won by 10 runs", "Sunrisers Hyderabad", "Royal Challengers Bangalore", "Rashid Khan ", "6", "5", "1", "0", "120.00"),
    @(" Abu Dhabi", " November 08 2020", "Capitals won by 17 runs", "Sunrisers Hyderabad", "Delhi Capitals", "Rashid Khan ", "11", "7", "1", "1", "157.14"),
    @(" Dubai (DSC)", " October 24 2020", "Kings XI won by 12 runs", "Sunrisers Hyderabad", "Kings XI Punjab", "Rashid Khan ", "0", "1", "0", "0", "0.00"),
    @(" Dubai (DSC)", " October 08 2020", "Sunrisers won by 69 runs", "Sunrisers Hyderabad", "Kings XI Punjab", "Rashid Khan ", "0", "0", "0", "0", "-")
)

$startRow = 9
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $cell = $ws.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $values[$c]
        $cell.ClearFormats()
    }
}
